# Actualización de horarios Línea 141 - 194
# Actualiza las 3 hojas (LP1912, LP1912-215, 6203-6173) con el nuevo scrape
# de las 04:36:47, insertando las nuevas filas intercaladas por Hora_Llegada
# y refrescando los contadores de cabecera.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Hoja 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 04:36:47"
$ws1.Cells.Item(3,1).Value = "Total filas: 18"

$s1 = @(
    @("04:01:01","04:01","81_EL PELIGRO",0,"LP1912"),
    @("04:36:47","04:45","215A_EL PATO",9,"LP1912"),
    @("04:01:01","04:46","215A_EL PATO",45,"LP1912"),
    @("04:01:01","04:53","11_ETCHEVERRY",52,"LP1912"),
    @("04:36:47","05:14","14_ABASTO",38,"LP1912"),
    @("04:01:01","05:16","17_ROMERO",75,"LP1912"),
    @("04:01:01","05:22","23_HERNANDEZ",81,"LP1912"),
    @("04:36:47","05:34","215B_EL PATO",58,"LP1912"),
    @("04:01:01","05:35","215B_EL PATO",94,"LP1912"),
    @("04:01:01","05:41","14_ABASTO",100,"LP1912"),
    @("04:01:01","05:46","15_ABASTO",105,"LP1912"),
    @("04:36:47","06:04","16_SANTA ANA",88,"LP1912"),
    @("04:36:47","06:11","215A_EL PATO",95,"LP1912"),
    @("04:36:47","06:14","225_HARAS DEL SUR",98,"LP1912"),
    @("04:36:47","06:21","26_HERNANDEZ",105,"LP1912"),
    @("04:36:47","06:27","23_HERNANDEZ",111,"LP1912"),
    @("04:36:47","06:29","86_EST CHICA-ESC AGRARIA",113,"LP1912"),
    @("04:36:47","06:31","16_SANTA ANA",115,"LP1912")
)

$r = 6
foreach ($row in $s1) {
    $ws1.Cells.Item($r,1).Value = $row[0]
    $ws1.Cells.Item($r,2).Value = $row[1]
    $ws1.Cells.Item($r,3).Value = $row[2]
    $ws1.Cells.Item($r,4).Value = $row[3]
    $ws1.Cells.Item($r,5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Hoja 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 04:36:47"
$ws2.Cells.Item(3,1).Value = "Total filas: 5"

$s2 = @(
    @("04:36:47","04:45","215A_EL PATO",9,"LP1912"),
    @("04:01:01","04:46","215A_EL PATO",45,"LP1912"),
    @("04:36:47","05:34","215B_EL PATO",58,"LP1912"),
    @("04:01:01","05:35","215B_EL PATO",94,"LP1912"),
    @("04:36:47","06:11","215A_EL PATO",95,"LP1912")
)

$r = 6
foreach ($row in $s2) {
    $ws2.Cells.Item($r,1).Value = $row[0]
    $ws2.Cells.Item($r,2).Value = $row[1]
    $ws2.Cells.Item($r,3).Value = $row[2]
    $ws2.Cells.Item($r,4).Value = $row[3]
    $ws2.Cells.Item($r,5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Hoja 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 04:36:47"
$ws3.Cells.Item(3,1).Value = "Total filas: 4"

$s3 = @(
    @("04:36:47","05:43","215A_LA PLATA",67,"L6173"),
    @("04:01:01","05:44","215A_LA PLATA",103,"L6173"),
    @("04:36:47","06:10","215A_LA PLATA",94,"L6173"),
    @("04:36:47","06:32","215C_LA PLATA",116,"L6203")
)

$r = 6
foreach ($row in $s3) {
    $ws3.Cells.Item($r,1).Value = $row[0]
    $ws3.Cells.Item($r,2).Value = $row[1]
    $ws3.Cells.Item($r,3).Value = $row[2]
    $ws3.Cells.Item($r,4).Value = $row[3]
    $ws3.Cells.Item($r,5).Value = $row[4]
    $r = $r + 1
}
